$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "T=E/(CV*M)-273" label to use -273.3 instead of -273
$ws.Range("E9").Value = "T=E/(CV*M)-273.3"

# New input cells D2 / E2, and feed D4 from them
$ws.Range("D2").Value = 6000000
$ws.Range("E2").Value = 1
$ws.Range("D4").Formula = "=D2*E2"

# Update T0 (K4) from 293 to 293.3
$ws.Range("K4").Value = 293.3

# Update the -273 constant to -273.3 in the calc-temperature column formulas
$ws.Range("E10").Formula = "=B10/`$H`$4/D10-273.3"
$ws.Range("E11").Formula = "=B11/`$H`$4/D11-273.3"

# Widen column F (6th column) to fit the new D2/E2 inputs
$ws.Columns.Item(6).ColumnWidth = 15

# Update the active selection
[void]$ws.Range("G6").Select()
